$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 588 lost its "today-forward-fill" formulas in columns B, H, J and K -
# the author cleared those cells (they were evaluating to "" anyway because
# there was no data yet for that date), leaving plain empty cells behind
# while keeping every other cell / style in the row untouched.
$ws.Range("B588").ClearContents()
$ws.Range("H588").ClearContents()
$ws.Range("J588").ClearContents()
$ws.Range("K588").ClearContents()

# The view was scrolled down so the frozen pane's visible corner sits near
# the bottom of the data (around row 582) and the selection moved to P584.
$win = $excel.ActiveWindow
$win.ScrollRow = 582
$win.ScrollColumn = 2
[void]$ws.Range("P584").Select()
